$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# Header row (row 1) - extend from G to N with new column headers
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2 (index 95)
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "潘孟安"
$ws.Range("D2").Value = "第一銀行恆春分行屏東縣恆春鎮中正路"
$ws.Range("E2").Value = 2365760
$ws.Range("F2").Value = "98年12月24日"
$ws.Range("G2").Value = "房貸"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-26"
$ws.Range("K2").Value = "潘孟安"
$ws.Range("L2").Value = 1376
$ws.Range("M2").Value = "tmpf07c1"
$ws.Range("N2").Value = 95

# Row 3 (index 96)
$ws.Range("B3").Value = "房屋貸款"
$ws.Range("C3").Value = "潘孟安"
$ws.Range("D3").Value = "彰化銀行車城分行屏東縣車城鄉福興村中山路"
$ws.Range("E3").Value = 4634451
$ws.Range("F3").Value = "98年07月06日"
$ws.Range("G3").Value = "房貸"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2013-12-26"
$ws.Range("K3").Value = "潘孟安"
$ws.Range("L3").Value = 1376
$ws.Range("M3").Value = "tmpf07c1"
$ws.Range("N3").Value = 96
